$d = $word.ActiveDocument

# Locate the list paragraph whose text is exactly "jQuery" (the
# "Prrafodelista" bullet that currently reads "jQuery" - it is the
# paragraph we need to turn into two separate bullets: "JavaScript"
# and "JQuery").
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs.Item($i)
    $text = $para.Range.Text.TrimEnd([char]13)
    if ($text -eq "jQuery") {
        $target = $para
        break
    }
}

if ($target -eq $null) {
    throw "Could not find the 'jQuery' paragraph to edit."
}

$full  = $target.Range
$start = $full.Start
$end   = $full.End
$whole = $d.Range($start, $end)

# Replace that single paragraph with two list paragraphs (same
# pStyle/numPr/list formatting): one reading "JavaScript" (split
# across two runs, "JavaS" + "cript") and a new one reading "JQuery"
# (split across three runs, "J" + "Q" + "uery", the last one keeping
# the original run's rsid so formatting stays tied to the existing
# "uery" text).
$xmlFrag = @'
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p w:rsidR="003C45A1" w:rsidRPr="003C45A1" w:rsidRDefault="003C45A1" w:rsidP="003C45A1"><w:pPr><w:pStyle w:val="Prrafodelista"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:rPr><w:lang w:val="es-CO"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="es-CO"/></w:rPr><w:t>JavaS</w:t></w:r><w:r><w:rPr><w:lang w:val="es-CO"/></w:rPr><w:t>cript</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Prrafodelista"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:rPr><w:lang w:val="es-CO"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="es-CO"/></w:rPr><w:t>J</w:t></w:r><w:r><w:rPr><w:lang w:val="es-CO"/></w:rPr><w:t>Q</w:t></w:r><w:r w:rsidRPr="003C45A1"><w:rPr><w:lang w:val="es-CO"/></w:rPr><w:t>uery</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@

$whole.InsertXML($xmlFrag)
